# Rename "Requested quantity" headers on existing sheets
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item(1)
$wsMonthly = $wb.Worksheets.Item(2)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add new "PO Forecast" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used by the other sheets in this workbook
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Copy header formatting (bold, border, centered) from an existing header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Copy the date-number-format used on column A of the existing sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A67").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$rows = @(
    @(44934.99999999999, 17, -5.434380341712059, 39.9656540114187),
    @(44941.99999999999, 16, -6.520779616610104, 40.57352259142046),
    @(44962.99999999999, 16, -7.984410074319374, 38.56204618248353),
    @(44969.99999999999, 16, -8.502468710088205, 38.42974950256683),
    @(44976.99999999999, 16, -8.152367523649305, 37.59337965461236),
    @(44983.99999999999, 16, -6.18632902595757, 38.61158886814689),
    @(44990.99999999999, 16, -6.993684592086185, 38.42934128973907),
    @(44997.99999999999, 16, -7.140845569472248, 38.61628020441611),
    @(45004.99999999999, 16, -6.604712947035052, 36.98479628884068),
    @(45011.99999999999, 16, -7.119881978891114, 38.43886021624522),
    @(45018.99999999999, 15, -8.167184794333986, 38.96765513004877),
    @(45025.99999999999, 15, -6.116953818599085, 36.07047371764916),
    @(45039.99999999999, 15, -8.443861413977537, 38.19584721442749),
    @(45046.99999999999, 15, -8.229811827602493, 38.65809950701187),
    @(45053.99999999999, 15, -7.678154916375469, 39.37767174018398),
    @(45060.99999999999, 15, -7.798726829708551, 38.08124285746646),
    @(45067.99999999999, 15, -7.37122730525572, 36.4510705844502),
    @(45074.99999999999, 15, -7.574741835996788, 35.74306611544311),
    @(45081.99999999999, 15, -8.622954103357936, 38.14778195639796),
    @(45088.99999999999, 15, -8.843264220781728, 37.49604295845713),
    @(45095.99999999999, 14, -8.526181887768034, 35.0263765699148),
    @(45102.99999999999, 14, -7.831679839940293, 36.32508023447402),
    @(45116.99999999999, 14, -9.894280621010704, 35.08948563507204),
    @(45123.99999999999, 14, -8.050769776692624, 37.1295339789205),
    @(45130.99999999999, 14, -8.513868214141041, 36.64640228086034),
    @(45137.99999999999, 14, -9.957160571810835, 37.93129638437311),
    @(45144.99999999999, 14, -8.079176242675436, 36.55671712905379),
    @(45151.99999999999, 14, -8.817308929008792, 36.73244650322795),
    @(45158.99999999999, 14, -10.77009096133065, 36.6130108013421),
    @(45179.99999999999, 13, -9.453361955399117, 35.96142723578465),
    @(45186.99999999999, 13, -11.38411814118986, 35.57670639212287),
    @(45193.99999999999, 13, -10.28097704593204, 35.34968323572054),
    @(45200.99999999999, 13, -9.226151845246054, 36.07780704712167),
    @(45207.99999999999, 13, -9.97624937205315, 35.12105690998984),
    @(45214.99999999999, 13, -11.10546736220902, 35.12057887590142),
    @(45221.99999999999, 13, -9.315507103792958, 36.07965701821462),
    @(45228.99999999999, 13, -10.14657482846472, 35.77258994735435),
    @(45326.99999999999, 11, -10.79438324549854, 34.55978089248303),
    @(45333.99999999999, 11, -12.08855409969555, 32.59534241853095),
    @(45340.99999999999, 11, -10.80896579184318, 34.76541640391627),
    @(45347.99999999999, 11, -11.57184338810868, 34.99119204499355),
    @(45354.99999999999, 11, -11.99705045552406, 33.95478675294954),
    @(45361.99999999999, 11, -12.36064508247912, 33.87504324206084),
    @(45368.99999999999, 11, -11.7419355604552, 33.62822237364513),
    @(45375.99999999999, 11, -11.71979007518557, 34.34063352995157),
    @(45389.99999999999, 11, -11.01597441730929, 34.20392200529559),
    @(45396.99999999999, 11, -10.06341609593196, 33.70730768337608),
    @(45403.99999999999, 10, -12.3314611293572, 32.28794127681326),
    @(45410.99999999999, 10, -11.16293148580633, 32.64757298879076),
    @(45417.99999999999, 10, -12.78469887518626, 31.90086070830991),
    @(45424.99999999999, 10, -13.56188409918824, 34.12895106229818),
    @(45431.99999999999, 10, -12.0063189956817, 31.19357813244217),
    @(45438.99999999999, 10, -11.38662228633172, 35.03629253270461),
    @(45445.99999999999, 10, -12.97960951076195, 33.02272420448259),
    @(45452.99999999999, 10, -13.18600566284582, 33.69652481304414),
    @(45459.99999999999, 10, -14.22563942834579, 32.69659855809933),
    @(45480.99999999999, 9, -13.85319524320592, 33.44367147220205),
    @(45487.99999999999, 9, -13.81909433214855, 32.60901045323585),
    @(45494.99999999999, 9, -13.30226059924795, 31.96471394522214),
    @(45501.99999999999, 9, -11.54115584782464, 32.49007525401245),
    @(45508.99999999999, 9, -14.67783218852563, 31.5558965827442),
    @(45515.99999999999, 9, -12.9766168092479, 30.94951795252395),
    @(45522.99999999999, 9, -12.45612923431342, 30.34591257914668),
    @(45529.99999999999, 9, -14.11886100414988, 31.92096757501079),
    @(45536.99999999999, 9, -15.48021136540036, 32.61413084248758),
    @(45543.99999999999, 9, -11.81269903312265, 32.54427223229519)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $r = $i + 2
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}

# Restore the originally-active sheet/selection
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select()
